# Applies the numeric updates to the three-digit / one-digit division
# worksheet. Each original expression is unique within the document, so a
# straightforward Find & Replace (no wildcards) for each pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("574÷9=", "354÷4="),
    @("927÷9=", "445÷9="),
    @("622÷6=", "416÷6="),
    @("502÷5=", "862÷8="),
    @("968÷6=", "574÷8="),
    @("673÷6=", "501÷9="),
    @("157÷6=", "428÷4="),
    @("676÷4=", "927÷4="),
    @("390÷4=", "879÷8="),
    @("486÷4=", "825÷2="),
    @("338÷2=", "860÷3="),
    @("784÷3=", "969÷4="),
    @("234÷8=", "134÷3="),
    @("740÷4=", "501÷8="),
    @("188÷6=", "706÷5="),
    @("740÷2=", "144÷6="),
    @("942÷6=", "930÷9="),
    @("362÷6=", "142÷6="),
    @("717÷8=", "279÷2="),
    @("376÷9=", "869÷5="),
    @("478÷8=", "185÷7="),
    @("236÷5=", "275÷3="),
    @("617÷6=", "347÷6="),
    @("796÷3=", "225÷9="),
    @("439÷6=", "249÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
